$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.2590000629425049,
    0.14599990844726562,
    0.11400008201599121,
    0.10400009155273438,
    0.11599993705749512,
    0.09599995613098145,
    0.09400010108947754,
    0.08599996566772461,
    0.09299993515014648,
    0.08700013160705566,
    0.1119999885559082,
    0.08999991416931152,
    0.08599996566772461,
    0.08099985122680664,
    0.07500004768371582,
    0.10400009155273438,
    0.10400009155273438,
    0.15000009536743164,
    0.1510000228881836,
    0.10400009155273438,
    0.0839998722076416,
    0.09200000762939453,
    0.1119999885559082,
    0.0969998836517334,
    0.11500000953674316,
    0.08500003814697266,
    0.09999990463256836,
    0.07399988174438477,
    0.07500004768371582,
    0.09200000762939453,
    0.0970001220703125,
    0.09099984169006348,
    0.10000014305114746,
    0.08699989318847656,
    0.0839998722076416,
    0.08599996566772461,
    0.09800004959106445,
    0.0969998836517334,
    0.0989999771118164,
    0.0820000171661377,
    0.08699989318847656,
    0.10100007057189941,
    0.09199976921081543,
    0.0820000171661377,
    0.07999992370605469,
    0.09599995613098145,
    0.11100006103515625,
    0.11999988555908203,
    0.11999988555908203,
    0.11899995803833008,
    0.11500000953674316,
    0.10299992561340332,
    0.1119999885559082,
    0.09399986267089844,
    0.1359999179840088,
    0.0840001106262207,
    0.0970001220703125,
    0.08899998664855957,
    0.08700013160705566,
    0.08299994468688965,
    0.08000016212463379,
    0.09599995613098145,
    0.09299993515014648,
    0.0969998836517334,
    0.0840001106262207,
    0.08899998664855957,
    0.08500003814697266,
    0.09299993515014648,
    0.08799982070922852,
    0.08500003814697266,
    0.10500001907348633,
    0.08099985122680664,
    0.06999993324279785,
    0.08700013160705566,
    0.09599995613098145,
    0.0970001220703125,
    0.09399986267089844,
    0.10500001907348633,
    0.09099984169006348,
    0.08999991416931152
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

# Remove the now-unused rows 82:121 (previous data extended to row 121)
$ws.Range("A82:A121").ClearContents()
